# Regenerate orders with updated distance/sizes.
# The experiment's distance codes and one size code were renumbered:
#   D51 -> D55
#   D80 -> D86
#   D64 -> D69
#   S30 -> S31
# These tokens appear as substrings inside many cell values across the
# whole sheet (Condition, Filename_Left, Filename_Right, Distance, Size
# columns), so perform a global text replace over the used range, the
# same way Excel's Find & Replace (Ctrl+H) would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$rng.Replace("D51", "D55")
$rng.Replace("D80", "D86")
$rng.Replace("D64", "D69")
$rng.Replace("S30", "S31")
